# Commit: "Code changes to incorporate SOAP request handling"
#
# The TestSuite1 sheet had a numbering bug: rows 5 and 6 both carried the
# "TC04" test-case id. This edit renumbers every test case from row 6
# onward (TC05 .. TC84), which also requires three brand new ids
# (TC82, TC83, TC84) to be created for the last three rows. A stray
# formatting override on D83 is cleared back to the common style, and the
# two sheets' saved view state (frozen-pane scroll position / selection)
# is updated to where the author ended up after the edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestSuite1")
$ws2 = $wb.Worksheets.Item("TestData")

# --- Renumber the Test Case column (A6:A85): row N -> "TC" + (N-1) -------
for ($r = 6; $r -le 85; $r++) {
    $n = $r - 1
    $label = "TC{0:D2}" -f $n
    $ws1.Cells.Item($r, 1).Value = $label
}

# --- D83 picked up a one-off style (s="17"); reset it to the normal ------
# --- wrap-text style used by every other cell in that column (s="1"). ----
$ws1.Range("D84").Copy() | Out-Null
$ws1.Range("D83").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Saved view state: TestData's selection moved, TestSuite1 stayed the
# --- active sheet with its frozen pane scrolled to row 83 and B90 picked.
$ws2.Range("J16").Select() | Out-Null

$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 83
$ws1.Range("B90").Select() | Out-Null
